# Updated cryptos list on Mon Feb 12 04:48:56 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) cells for rows 2-51
# of the active worksheet to the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. A leading apostrophe is used so Excel
# stores values such as "38.70" or "0.140" as literal text (matching the
# source data) instead of silently coercing them into numbers and dropping
# trailing zeros.
$updates = [ordered]@{
    'D2' = '''48.050.39'
    'E2' = '''  -0.67%  '
    'D3' = '''2.495.50'
    'E3' = '''  -1.16%  '
    'E4' = '''  -0.07%  '
    'D5' = '''319.77'
    'E5' = '''  -1.35%  '
    'D6' = '''105.72'
    'E6' = '''  -3.45%  '
    'E7' = '''  -1.11%  '
    'E8' = '''  -0.01%  '
    'E9' = '''  -4.53%  '
    'D10' = '''38.70'
    'D11' = '''20.06'
    'E11' = '''  +1.04%  '
    'E12' = '''  -2.09%  '
    'E13' = '''  -0.44%  '
    'E14' = '''  -2.38%  '
    'D15' = '''2.887.79'
    'E15' = '''  -1.12%  '
    'D16' = '''2.505.96'
    'E16' = '''  -0.72%  '
    'E17' = '''  -3.00%  '
    'D18' = '''47.898.49'
    'E18' = '''  -0.70%  '
    'D19' = '''13.10'
    'E19' = '''  -2.39%  '
    'E20' = '''  +7.94%  '
    'D21' = '''6.65'
    'E21' = '''  -0.32%  '
    'E22' = '''  -1.42%  '
    'D23' = '''71.14'
    'E23' = '''  -2.29%  '
    'D24' = '''272.15'
    'E24' = '''  -0.32%  '
    'E25' = '''  -2.81%  '
    'E26' = '''  -0.13%  '
    'E27' = '''  -1.57%  '
    'D28' = '''2.29'
    'E28' = '''  +0.92%  '
    'E29' = '''  -4.63%  '
    'D30' = '''0.140'
    'E30' = '''  -4.30%  '
    'D31' = '''34.71'
    'E31' = '''  -1.88%  '
    'D32' = '''49.11'
    'E32' = '''  -1.41%  '
    'E33' = '''  -0.06%  '
    'D34' = '''19.12'
    'E34' = '''  -4.44%  '
    'E35' = '''  -2.64%  '
    'D36' = '''0.0774'
    'E36' = '''  -2.21%  '
    'E37' = '''  -2.79%  '
    'D38' = '''4.56'
    'E38' = '''  -3.75%  '
    'E39' = '''  -4.60%  '
    'D40' = '''122.33'
    'E40' = '''  +2.71%  '
    'E41' = '''  -2.21%  '
    'D42' = '''22.19'
    'E42' = '''  -0.15%  '
    'E43' = '''  +1.01%  '
    'E44' = '''  +0.89%  '
    'D45' = '''1.997.52'
    'E45' = '''  -0.32%  '
    'D46' = '''3.14'
    'E46' = '''  +0.32%  '
    'D47' = '''1.89'
    'E47' = '''  +0.18%  '
    'E48' = '''  -1.12%  '
    'E49' = '''  -2.38%  '
    'D50' = '''5.18'
    'E50' = '''  -1.91%  '
    'D51' = '''78.76'
    'E51' = '''  -2.76%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

